$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Update the title text
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Slide 1 v.2"

# Clear the subtitle text (becomes an empty paragraph)
$s.Shapes.Item(2).TextFrame.TextRange.Text = ""
